$wb = $excel.ActiveWorkbook

# Update the header labels on the "flow" sheet to the new shared-string
# values (chb_irr1, chb_con2, chs_irr3, chs_con4, chs_irr5, chs_con6).
$flow = $wb.Worksheets.Item("flow")
$flow.Range("B1").Value = "chb_irr1"
$flow.Range("C1").Value = "chb_con2"
$flow.Range("D1").Value = "chs_irr3"
$flow.Range("E1").Value = "chs_con4"
$flow.Range("F1").Value = "chs_irr5"
$flow.Range("G1").Value = "chs_con6"

# Move the active tab / selection from "par" to "flow".
$flow.Activate()
$flow.Range("G3").Select()
